$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.191.74'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').Value = '1.659.78'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.73'
$ws.Range('E5').Value = '  -1.41%  '
$ws.Range('E6').Value = '  -2.89%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2645'
$ws.Range('E8').Value = '  -1.43%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06271'
$ws.Range('E9').Value = '  -1.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.82'
$ws.Range('E10').Value = '  -4.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07755'
$ws.Range('E11').Value = '  -0.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.464'
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('D13').Value = '1.652.92'
$ws.Range('E13').Value = '  -1.29%  '
$ws.Range('D14').Value = '1.885.62'
$ws.Range('E14').Value = '  -1.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5438'
$ws.Range('E15').Value = '  -2.53%  '
$ws.Range('D16').Value = '0.0₅8102'
$ws.Range('E16').Value = '  -2.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.89'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').Value = '26.193.50'
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('E19').Value = '  +0.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.622'
$ws.Range('E20').Value = '  -3.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '192.04'
$ws.Range('E21').Value = '  -0.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.08'
$ws.Range('E22').Value = '  -2.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.007'
$ws.Range('E23').Value = '  -4.91%  '
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '140.03'
$ws.Range('E25').Value = '  +1.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1220'
$ws.Range('E26').Value = '  -4.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.252'
$ws.Range('E27').Value = '  -2.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.16'
$ws.Range('E28').Value = '  -1.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.433'
$ws.Range('E29').Value = '  +0.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05961'
$ws.Range('E30').Value = '  -5.07%  '
$ws.Range('E31').Value = '  -1.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.568'
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.261'
$ws.Range('E33').Value = '  -4.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.595'
$ws.Range('E34').Value = '  -5.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9659'
$ws.Range('E35').Value = '  -4.46%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.424'
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.768'
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5696'
$ws.Range('E38').Value = '  -8.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01592'
$ws.Range('E39').Value = '  -1.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.984'
$ws.Range('E40').Value = '  -1.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8589'
$ws.Range('E41').Value = '  -0.32%  '
$ws.Range('E42').Value = '  +0.22%  '
$ws.Range('D43').Value = '1.010.91'
$ws.Range('E43').Value = '  -7.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.30'
$ws.Range('E44').Value = '  -0.31%  '
$ws.Range('D45').Value = '1.799.88'
$ws.Range('E45').Value = '  -1.19%  '
$ws.Range('D46').Value = '0.0₈108'
$ws.Range('E46').Value = '  -0.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.68'
$ws.Range('E47').Value = '  -3.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.006'
$ws.Range('E48').Value = '  +0.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.041'
$ws.Range('E49').Value = '  -2.07%  '
$ws.Range('E50').Value = '  -0.53%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.457'
$ws.Range('E51').Value = '  -5.03%  '

Write-Output "Applied 94 cell updates"
